$d = $word.ActiveDocument

function Get-SubRange($ctx, $target) {
    $r = $d.Content
    $r.Find.Execute($ctx, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    $s = $r.Start
    $idx = $ctx.IndexOf($target)
    $subStart = $s + $idx
    $subEnd = $subStart + $target.Length
    return $d.Range($subStart, $subEnd)
}

# ---------------------------------------------------------------------------
# 1) "What is LOREM IPSUM[ changed]?" heading - the blue " changed" run loses
#    its bold formatting (Word also stamps an explicit strike-off on the run
#    once its direct character formatting is touched).
# ---------------------------------------------------------------------------
$r = Get-SubRange "What is LOREM IPSUM changed?" " changed"
$r.Font.Bold = 0
$r.Font.StrikeThrough = 0

# ---------------------------------------------------------------------------
# 2) "...electronic dasd asd typesetting, remaining essentially unchanged..."
#    -> "asd", the following space and "typesetting, remaining " get struck
#    through (the rest of the sentence keeps its original formatting, which
#    is why Word splits the run right after "remaining ").
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("asd typesetting, remaining ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$r.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 3) The five blue "changed here too" / "changed here" runs: strikethrough is
#    explicitly pinned off on each of them.
# ---------------------------------------------------------------------------
$r = Get-SubRange "Letraset sheets containing Lorem Ipsum changed here too passages" "changed here too "
$r.Font.StrikeThrough = 0

$r = Get-SubRange "Aldus PageMaker including versions of Lorem Ipsum changed here." " changed here"
$r.Font.StrikeThrough = 0

$r = Get-SubRange "as their default model changed here too text" "changed here too "
$r.Font.StrikeThrough = 0

$r = Get-SubRange "obscure Latin words changed here too, consectetur" " changed here too"
$r.Font.StrikeThrough = 0

$r = Get-SubRange "de Finibus Bonorum et changed here too Malorum" "changed here too "
$r.Font.StrikeThrough = 0
